# "Shooting and smoke dash"
# Append a new log row (Monday 2025-05-12, 08:30-10:30) describing the
# "shooting & smoke dash" work, mirroring the formatting of the row above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14
$prevRow = 13

# --- Day -----------------------------------------------------------------
$ws.Cells.Item($row, 1).Value = "Monday"

# --- Date ------------------------------------------------------------------
$ws.Cells.Item($prevRow, 2).Copy()
$ws.Cells.Item($row, 2).PasteSpecial(-4122)
$ws.Cells.Item($row, 2).Value = 45789

# --- From --------------------------------------------------------------
$ws.Cells.Item($prevRow, 3).Copy()
$ws.Cells.Item($row, 3).PasteSpecial(-4122)
$ws.Cells.Item($row, 3).Value = 0.35416666666666669

# --- Until -------------------------------------------------------------
$ws.Cells.Item($prevRow, 4).Copy()
$ws.Cells.Item($row, 4).PasteSpecial(-4122)
$ws.Cells.Item($row, 4).Value = 0.4375

# --- Time spent (=Until-From, same style as the rest of the column) ------
$ws.Cells.Item($prevRow, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)
$ws.Cells.Item($row, 5).Formula = "=D$row-C$row"

# --- Realised ----------------------------------------------------------
$ws.Cells.Item($row, 6).Value = "Working on shooting & smoke dash"

# Tidy up the clipboard / marching ants left by the copy operations and
# move the selection to the next empty row, like a user would after typing.
$excel.CutCopyMode = 0
[void]$ws.Range("F15").Select()

Write-Output "Added log entry for 2025-05-12 (shooting & smoke dash)"
